# "Generate Report for Handback"
#
# The nightly localization-status report now records the handback
# (target-file-returned) leg for each locale: a "Latest Target File" /
# "Latest Handback File" column pair is populated on the zh-cn and de-de
# sheets, the per-row Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and the handback timestamp columns
# get real datetimes instead of the 0001-01-01 sentinel.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkUrl($worksheet, $cellRef) {
    $targetAddr = $worksheet.Range($cellRef).Address()
    foreach ($h in $worksheet.Hyperlinks) {
        if ($h.Range.Address() -eq $targetAddr) {
            return $h.Address
        }
    }
    return $null
}

function Add-MatchingHyperlink($worksheet, $destCellRef, $srcCellRef) {
    $url = Get-HyperlinkUrl $worksheet $srcCellRef
    $text = $worksheet.Range($srcCellRef).Text
    $worksheet.Hyperlinks.Add($worksheet.Range($destCellRef), $url, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Status: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (every cell that showed the old status, across all three sheets)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback
#    File" (G) for both data rows, with the same hyperlink target/text
#    as the matching source-file / target-file columns, and set the
#    "Latest Handback DateTime" (H) now that handback has happened.
# ---------------------------------------------------------------------
$wsZhCn.Range("F2").Value = $wsZhCn.Range("A2").Text
$wsZhCn.Range("G2").Value = $wsZhCn.Range("D2").Text
$wsZhCn.Range("F3").Value = $wsZhCn.Range("A3").Text
$wsZhCn.Range("G3").Value = $wsZhCn.Range("D3").Text

Add-MatchingHyperlink $wsZhCn "F2" "A2"
Add-MatchingHyperlink $wsZhCn "G2" "D2"
Add-MatchingHyperlink $wsZhCn "F3" "A3"
Add-MatchingHyperlink $wsZhCn "G3" "D3"

$zhCnLinkRange = $wsZhCn.Range("F2:G3")
$zhCnLinkRange.Font.Underline = 2
$zhCnLinkRange.Font.Color = 15570276

$wsZhCn.Range("H2").Value = "2016-03-12 02:10:04"
$wsZhCn.Range("H3").Value = "2016-03-12 02:10:04"

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape of change, with its own handback time.
# ---------------------------------------------------------------------
$wsDeDe.Range("F2").Value = $wsDeDe.Range("A2").Text
$wsDeDe.Range("G2").Value = $wsDeDe.Range("D2").Text
$wsDeDe.Range("F3").Value = $wsDeDe.Range("A3").Text
$wsDeDe.Range("G3").Value = $wsDeDe.Range("D3").Text

Add-MatchingHyperlink $wsDeDe "F2" "A2"
Add-MatchingHyperlink $wsDeDe "G2" "D2"
Add-MatchingHyperlink $wsDeDe "F3" "A3"
Add-MatchingHyperlink $wsDeDe "G3" "D3"

$deDeLinkRange = $wsDeDe.Range("F2:G3")
$deDeLinkRange.Font.Underline = 2
$deDeLinkRange.Font.Color = 15570276

$wsDeDe.Range("H2").Value = "2016-03-12 02:10:09"
$wsDeDe.Range("H3").Value = "2016-03-12 02:10:09"

Write-Host "Handback report generated."
